$wb = $excel.ActiveWorkbook

# --- "Atlas" sheet: flip a handful of True/False flags in the
# "Binned graph" (D) / "Ratio graph" (E) columns. These are stored as plain
# text (shared strings), so copy/paste-values from an existing "False" cell
# instead of assigning the literal which Excel would auto-coerce to a
# genuine Boolean.
$atlas = $wb.Worksheets.Item("Atlas")

$falseSource = $atlas.Range("E2")
$falseSource.Copy()
$atlas.Range("D2").PasteSpecial(-4163)
$atlas.Range("D3").PasteSpecial(-4163)
$atlas.Range("E5").PasteSpecial(-4163)
$atlas.Range("E6").PasteSpecial(-4163)
$atlas.Range("E23").PasteSpecial(-4163)

# C6 loses its one-off "applyFont" style variant, matching the plain
# centered style used by the rest of column C.
$atlas.Range("C2").Copy()
$atlas.Range("C6").PasteSpecial(-4122)  # xlPasteFormats

# Leave the selection the way it was left on the Atlas sheet (D2:D4, active
# cell D4).
$atlas.Range("D2:D4").Select()

# --- "Excel" sheet: F22's style loses its explicit font, matching the
# plain centered style used elsewhere (e.g. A22:E22).
$excelSheet = $wb.Worksheets.Item("Excel")
$excelSheet.Range("F22").Font.Bold = $false
